# VerveStacks BGR model update - 2025-08-26 23:35
# Re-shuffles the day/night timeslice orderings on "ev_charging_uc" and the
# seasonal RE-profile rows on "re_profiles".

$wb = $excel.ActiveWorkbook

# --- Sheet "ev_charging_uc": update the two comma-separated timeslice
#     sequences that feed the G7 (=C14) / G8 (=C13) formulas.
$wsUC = $wb.Worksheets.Item("ev_charging_uc")
$wsUC.Range("C13").Value = "WaP,WaD,RaP,SaD,RaD,FaD,FaP,SaP"
$wsUC.Range("C14").Value = "RaN,FaP,SaP,RaP,FaN,SaN,WaN,WaP"

# --- Sheet "re_profiles": the M:N (season / share) rows were re-ordered.
$wsRE = $wb.Worksheets.Item("re_profiles")
$wsRE.Range("M4").Value = "W"
$wsRE.Range("N4").Value = 0.27238459437312212
$wsRE.Range("M5").Value = "S"
$wsRE.Range("N5").Value = 0.34121824638077031
$wsRE.Range("M6").Value = "F"
$wsRE.Range("N6").Value = 0.1858508604206501
$wsRE.Range("M7").Value = "R"
$wsRE.Range("N7").Value = 0.4005462988254575
